# Weekly fruit/vegetable price update: insert a new record as row 267,
# pushing the existing rows 267:379 down to 268:380.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 267 (shifts rows 267-379 down to 268-380)
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row 267 with this week's record.
$row = 267
$ws.Cells.Item($row, 1).Value  = 10
$ws.Cells.Item($row, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item($row, 3).Value  = 'La Araucanía'
$ws.Cells.Item($row, 4).Value  = 44839
$ws.Cells.Item($row, 5).Value  = 9
$ws.Cells.Item($row, 6).Value  = 100112017
$ws.Cells.Item($row, 7).Value  = 'Apio'
$ws.Cells.Item($row, 8).Value  = 'Americana (o)'
$ws.Cells.Item($row, 9).Value  = 'Primera'
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 12000
$ws.Cells.Item($row, 12).Value = 12000
$ws.Cells.Item($row, 13).Value = 12000
$ws.Cells.Item($row, 14).Value = '$/docena de matas'
$ws.Cells.Item($row, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item($row, 16).Value = 2000
$ws.Cells.Item($row, 17).Value = 6
$ws.Cells.Item($row, 18).Value = 'Hortaliza'
